$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B13 value (350 -> 253)
$ws.Range("B13").Value = "253"

# Add new rows 14-21 with invalid test data
$ws.Range("A14").Value = "invalid.productcode"
$ws.Range("B14").Value = "test,asd2314, ,00,@#`$123,`$#@asd,1234"

$ws.Range("A15").Value = "invalid.flag"
$ws.Range("B15").Value = "123,aw234,@#!w34,test, ,012@!#,-2098"

$ws.Range("A16").Value = "invalid.integervalue"
$ws.Range("B16").Value = "@-123,aw234,@#!w34,test, ,012@!#,-(*2098"

$ws.Range("A17").Value = "invalid.uomid"
$ws.Range("B17").Value = "test,asd2314, ,00,@#`$123,`$#@asd,1234"

$ws.Range("A18").Value = "invalid.godownflag"
$ws.Range("B18").Value = "@-123,aw234,@#!w34,test, ,012@!#,-(*2098"

$ws.Range("A19").Value = "invalid.qty"
$ws.Range("B19").Value = "test,asd2314, ,00,@#`$123,`$#@asd,1234"

$ws.Range("A20").Value = "invalid.rate"
$ws.Range("B20").Value = "@-123,aw234,@#!w34,test, ,012@!#,-(*2098"

$ws.Range("A21").Value = "invalid.categorycode"
$ws.Range("B21").Value = "123,aw234,@#!w34,test, ,012@!#,-2098"

# Update selection to C13
$ws.Range("C13").Select()
